$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (rows 2-13) to remove stale shared strings / values
$ws.Range("A2:T13").Clear()

# Write the string columns first, in column-major order (A for all rows, then B, then C, then D)
# so that the shared-strings table is rebuilt in the same reference order as the target file.
# Column A
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("A7").Value = "Neutrophils"
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("A13").Value = "Resolving-Mac"

# Column B
$ws.Range("B2").Value = "Wnt8a"
$ws.Range("B3").Value = "Wnt8a"
$ws.Range("B4").Value = "Wnt8a"
$ws.Range("B5").Value = "Wnt8a"
$ws.Range("B6").Value = "Wnt8a"
$ws.Range("B7").Value = "Wnt8a"
$ws.Range("B8").Value = "Wnt8a"
$ws.Range("B9").Value = "Wnt8a"
$ws.Range("B10").Value = "Wnt8a"
$ws.Range("B11").Value = "Wnt8a"
$ws.Range("B12").Value = "Wnt8a"
$ws.Range("B13").Value = "Wnt8a"

# Column C
$ws.Range("C2").Value = "Fzd4"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("C12").Value = "Fzd4"
$ws.Range("C13").Value = "Fzd4"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("D10").Value = "ECs"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("D13").Value = "MuSCs"

# Write the numeric columns
# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09452133333333333
$ws.Range("H2").Value = 0.283564
$ws.Range("I2").Value = 0.005525598587839362
$ws.Range("J2").Value = 0.005525598587839361
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 30.801072
$ws.Range("N2").Value = 61.602144
$ws.Range("O2").Value = 0.5373480691764108
$ws.Range("P2").Value = 0.485871843331092
$ws.Range("Q2").Value = 2.911358393536
$ws.Range("R2").Value = 17.468150361216
$ws.Range("S2").Value = 0.002969169732219383
$ws.Range("T2").Value = 0.002684732771381189

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09452133333333333
$ws.Range("H3").Value = 0.283564
$ws.Range("I3").Value = 0.005525598587839362
$ws.Range("J3").Value = 0.005525598587839361
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.015213
$ws.Range("N3").Value = 36.045639
$ws.Range("O3").Value = 0.2096145064786482
$ws.Range("P3").Value = 0.2843011610923331
$ws.Range("Q3").Value = 1.135693953044
$ws.Range("R3").Value = 10.221245577396
$ws.Range("S3").Value = 0.001158245620989064
$ws.Range("T3").Value = 0.001570934094252886

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09452133333333333
$ws.Range("H4").Value = 0.283564
$ws.Range("I4").Value = 0.005525598587839362
$ws.Range("J4").Value = 0.005525598587839361
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1305583333333333
$ws.Range("N4").Value = 0.391675
$ws.Range("O4").Value = 0.002277689176907768
$ws.Range("P4").Value = 0.003089240761436898
$ws.Range("Q4").Value = 0.01234054774444444
$ws.Range("R4").Value = 0.1110649297
$ws.Range("S4").Value = [double]"1.258559609945856E-05"
$ws.Range("T4").Value = [double]"1.706990438889152E-05"

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09452133333333333
$ws.Range("H5").Value = 0.283564
$ws.Range("I5").Value = 0.005525598587839362
$ws.Range("J5").Value = 0.005525598587839361
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.373679
$ws.Range("N5").Value = 28.747358
$ws.Range("O5").Value = 0.2507597351680332
$ws.Range("P5").Value = 0.2267377548151379
$ws.Range("Q5").Value = 1.358619303985333
$ws.Range("R5").Value = 8.151715823911999
$ws.Range("S5").Value = 0.001385597638531457
$ws.Range("T5").Value = 0.001252861817816393

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.98384066666667
$ws.Range("H6").Value = 50.951522
$ws.Range("I6").Value = 0.9928540224128104
$ws.Range("J6").Value = 0.9928540224128104
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 30.801072
$ws.Range("N6").Value = 61.602144
$ws.Range("O6").Value = 0.5373480691764108
$ws.Range("P6").Value = 0.485871843331092
$ws.Range("Q6").Value = 523.1204992105279
$ws.Range("R6").Value = 3138.722995263167
$ws.Range("S6").Value = 0.5335081919175565
$ws.Range("T6").Value = 0.4823998140284015

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.98384066666667
$ws.Range("H7").Value = 50.951522
$ws.Range("I7").Value = 0.9928540224128104
$ws.Range("J7").Value = 0.9928540224128104
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.015213
$ws.Range("N7").Value = 36.045639
$ws.Range("O7").Value = 0.2096145064786482
$ws.Range("P7").Value = 0.2843011610923331
$ws.Range("Q7").Value = 204.064463168062
$ws.Range("R7").Value = 1836.580168512558
$ws.Range("S7").Value = 0.208116605913402
$ws.Range("T7").Value = 0.2822695513671553

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.98384066666667
$ws.Range("H8").Value = 50.951522
$ws.Range("I8").Value = 0.9928540224128104
$ws.Range("J8").Value = 0.9928540224128104
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1305583333333333
$ws.Range("N8").Value = 0.391675
$ws.Range("O8").Value = 0.002277689176907768
$ws.Range("P8").Value = 0.003089240761436898
$ws.Range("Q8").Value = 2.217381931038889
$ws.Range("R8").Value = 19.95643737935
$ws.Range("S8").Value = 0.002261412861099001
$ws.Range("T8").Value = 0.003067165116194238

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.98384066666667
$ws.Range("H9").Value = 50.951522
$ws.Range("I9").Value = 0.9928540224128104
$ws.Range("J9").Value = 0.9928540224128104
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.373679
$ws.Range("N9").Value = 28.747358
$ws.Range("O9").Value = 0.2507597351680332
$ws.Range("P9").Value = 0.2267377548151379
$ws.Range("Q9").Value = 244.1202739298127
$ws.Range("R9").Value = 1464.721643578876
$ws.Range("S9").Value = 0.2489678117207529
$ws.Range("T9").Value = 0.2251174919010592

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.02771833333333333
$ws.Range("H10").Value = 0.08315500000000001
$ws.Range("I10").Value = 0.001620378999350348
$ws.Range("J10").Value = 0.001620378999350348
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 30.801072
$ws.Range("N10").Value = 61.602144
$ws.Range("O10").Value = 0.5373480691764108
$ws.Range("P10").Value = 0.485871843331092
$ws.Range("Q10").Value = 0.85375438072
$ws.Range("R10").Value = 5.12252628432
$ws.Range("S10").Value = 0.0008707075266349141
$ws.Range("T10").Value = 0.0007872965313093441

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.02771833333333333
$ws.Range("H11").Value = 0.08315500000000001
$ws.Range("I11").Value = 0.001620378999350348
$ws.Range("J11").Value = 0.001620378999350348
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 12.015213
$ws.Range("N11").Value = 36.045639
$ws.Range("O11").Value = 0.2096145064786482
$ws.Range("P11").Value = 0.2843011610923331
$ws.Range("Q11").Value = 0.333041679005
$ws.Range("R11").Value = 2.997375111045
$ws.Range("S11").Value = 0.0003396549442571891
$ws.Range("T11").Value = 0.0004606756309249369

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.02771833333333333
$ws.Range("H12").Value = 0.08315500000000001
$ws.Range("I12").Value = 0.001620378999350348
$ws.Range("J12").Value = 0.001620378999350348
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1305583333333333
$ws.Range("N12").Value = 0.391675
$ws.Range("O12").Value = 0.002277689176907768
$ws.Range("P12").Value = 0.003089240761436898
$ws.Range("Q12").Value = 0.003618859402777778
$ws.Range("R12").Value = 0.032569734625
$ws.Range("S12").Value = [double]"3.690719709308927E-06"
$ws.Range("T12").Value = [double]"5.00574085376943E-06"

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.02771833333333333
$ws.Range("H13").Value = 0.08315500000000001
$ws.Range("I13").Value = 0.001620378999350348
$ws.Range("J13").Value = 0.001620378999350348
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 14.373679
$ws.Range("N13").Value = 28.747358
$ws.Range("O13").Value = 0.2507597351680332
$ws.Range("P13").Value = 0.2267377548151379
$ws.Range("Q13").Value = 0.3984144257483334
$ws.Range("R13").Value = 2.39048655449
$ws.Range("S13").Value = 0.000406325808748936
$ws.Range("T13").Value = 0.0003674010962622977
